$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate professor row (row 48: Silva, Jorge - duplicate of row 35)
# This shifts all subsequent rows up by one.
$ws.Rows.Item(48).Delete()

# Update the selection/view to match the post-edit state
$ws.Application.ActiveWindow.ScrollRow = 15
$ws.Range("B48:AJ53").Select()
